$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 350.42
$ws.Range("I15").Value = 350.42
$ws.Range("K15").Value = 1051.26
$ws.Range("M15").Value = -882.26

$ws.Range("H42").Value = 107
$ws.Range("I42").Value = 107
$ws.Range("K42").Value = 321
$ws.Range("M42").Value = -91

$ws.Range("H98").Value = 1188.3334
$ws.Range("I98").Value = 1226.05
$ws.Range("K98").Value = 1226.05
$ws.Range("M98").Value = 271.95

$ws.Range("H106").Value = 63495004
$ws.Range("I106").Value = 33336510
$ws.Range("K106").Value = 33336510
$ws.Range("M106").Value = -33335879

$ws.Range("H122").Value = 1188.3334
$ws.Range("I122").Value = 1226.05
$ws.Range("K122").Value = 3678.15
$ws.Range("M122").Value = -1228.15

$ws.Range("H135").Value = 1354.1132
$ws.Range("I135").Value = 1180.7805
$ws.Range("J135").Value = 1946.3334
$ws.Range("K135").Value = 10627.0245
$ws.Range("L135").Value = 17517.0006
$ws.Range("M135").Value = -8092.024500000001
$ws.Range("N135").Value = -22587.0006

$ws.Range("H137").Value = 1147.2699
$ws.Range("I137").Value = 966.587
$ws.Range("J137").Value = 1636.1765
$ws.Range("K137").Value = 2899.761
$ws.Range("L137").Value = 4908.529500000001
$ws.Range("M137").Value = -349.761
$ws.Range("N137").Value = -10008.5295

$ws.Range("H138").Value = 2484.261
$ws.Range("I138").Value = 970.93335
$ws.Range("J138").Value = 5321.75
$ws.Range("K138").Value = 2912.80005
$ws.Range("L138").Value = 15965.25
$ws.Range("M138").Value = 2227.19995
$ws.Range("N138").Value = -26245.25

$ws.Range("H141").Value = 1346.9464
$ws.Range("I141").Value = 915.3261
$ws.Range("J141").Value = 3332.4
$ws.Range("K141").Value = 2745.9783
$ws.Range("L141").Value = 9997.200000000001
$ws.Range("M141").Value = 2434.0217
$ws.Range("N141").Value = -20357.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2852.6785
$ws.Range("I61").Value = 2912.745
$ws.Range("J61").Value = 2240
$ws.Range("K61").Value = 2912.745
$ws.Range("L61").Value = 2240
$ws.Range("M61").Value = -2700.745
$ws.Range("N61").Value = -2664

$ws.Range("H74").Value = 1126.2325
$ws.Range("I74").Value = 1100.64
$ws.Range("J74").Value = 1161.7778
$ws.Range("K74").Value = 1100.64
$ws.Range("L74").Value = 1161.7778
$ws.Range("M74").Value = -226.6400000000001
$ws.Range("N74").Value = -2909.7778

$ws.Range("H77").Value = 1126.2325
$ws.Range("I77").Value = 1100.64
$ws.Range("J77").Value = 1161.7778
$ws.Range("K77").Value = 5503.200000000001
$ws.Range("L77").Value = 5808.889
$ws.Range("M77").Value = -1135.200000000001
$ws.Range("N77").Value = -14544.889

$ws.Range("H122").Value = 1976582.1
$ws.Range("I122").Value = 2568687.2
$ws.Range("J122").Value = 2898.3333
$ws.Range("K122").Value = 7706061.600000001
$ws.Range("L122").Value = 8694.999899999999
$ws.Range("M122").Value = -7703611.600000001
$ws.Range("N122").Value = -13594.9999

$ws.Range("H132").Value = 2502840.5
$ws.Range("I132").Value = 2426.4783
$ws.Range("J132").Value = 5885753.5
$ws.Range("K132").Value = 7279.4349
$ws.Range("L132").Value = 17657260.5
$ws.Range("M132").Value = -4749.4349
$ws.Range("N132").Value = -17662320.5

$ws.Range("H136").Value = 2852.6785
$ws.Range("I136").Value = 2912.745
$ws.Range("J136").Value = 2240
$ws.Range("K136").Value = 8738.235000000001
$ws.Range("L136").Value = 6720
$ws.Range("M136").Value = -6188.235000000001
$ws.Range("N136").Value = -11820

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 166668480
$ws.Range("I99").Value = 250001100
$ws.Range("K99").Value = 250001100
$ws.Range("M99").Value = -249999602

$ws.Range("H134").Value = 3678.5334
$ws.Range("I134").Value = 4072.647
$ws.Range("K134").Value = 12217.941
$ws.Range("M134").Value = -9682.940999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 183773.83
$ws.Range("I31").Value = 1407.017
$ws.Range("J31").Value = 582279.0600000001
$ws.Range("K31").Value = 1407.017
$ws.Range("L31").Value = 582279.0600000001
$ws.Range("M31").Value = -1112.017
$ws.Range("N31").Value = -582869.0600000001

$ws.Range("H34").Value = 183773.83
$ws.Range("I34").Value = 1407.017
$ws.Range("J34").Value = 582279.0600000001
$ws.Range("K34").Value = 1407.017
$ws.Range("L34").Value = 582279.0600000001
$ws.Range("M34").Value = -1205.017
$ws.Range("N34").Value = -582683.0600000001

$ws.Range("H58").Value = 1008.2656
$ws.Range("I58").Value = 636
$ws.Range("K58").Value = 636
$ws.Range("M58").Value = -433

$ws.Range("H107").Value = 15873829
$ws.Range("I107").Value = 23810028
$ws.Range("J107").Value = 1429.7142
$ws.Range("K107").Value = 23810028
$ws.Range("L107").Value = 1429.7142
$ws.Range("M107").Value = -23808108
$ws.Range("N107").Value = -5269.7142

$ws.Range("H132").Value = 1702.836
$ws.Range("I132").Value = 1453.5555
$ws.Range("J132").Value = 2403.9375
$ws.Range("K132").Value = 4360.666499999999
$ws.Range("L132").Value = 7211.8125
$ws.Range("M132").Value = -1830.666499999999
$ws.Range("N132").Value = -12271.8125

$ws.Range("H134").Value = 2393.422
$ws.Range("I134").Value = 3306.652
$ws.Range("J134").Value = 1438.6818
$ws.Range("K134").Value = 9919.956
$ws.Range("L134").Value = 4316.0454
$ws.Range("M134").Value = -7384.956
$ws.Range("N134").Value = -9386.045399999999

$ws.Range("H136").Value = 1008.2656
$ws.Range("I136").Value = 636
$ws.Range("K136").Value = 1908
$ws.Range("M136").Value = 642

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 1310.4
$ws.Range("I69").Value = 517.3333
$ws.Range("J69").Value = 2500
$ws.Range("K69").Value = 1551.9999
$ws.Range("L69").Value = 7500
$ws.Range("M69").Value = -740.9999
$ws.Range("N69").Value = -9122

$ws.Range("H72").Value = 1310.4
$ws.Range("I72").Value = 517.3333
$ws.Range("J72").Value = 2500
$ws.Range("K72").Value = 4655.9997
$ws.Range("L72").Value = 22500
$ws.Range("M72").Value = -599.9997000000003
$ws.Range("N72").Value = -30612

$ws.Range("H110").Value = 3790.9092
$ws.Range("J110").Value = 5166.6665
$ws.Range("L110").Value = 15499.9995
$ws.Range("N110").Value = -23679.9995

$ws.Range("H140").Value = 1993.4348
$ws.Range("I140").Value = 1993.4348
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 5980.3044
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -800.3044
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2819988.8
$ws.Range("I122").Value = 3814425.8
$ws.Range("J122").Value = 2417.1667
$ws.Range("K122").Value = 11443277.4
$ws.Range("L122").Value = 7251.500100000001
$ws.Range("M122").Value = -11440827.4
$ws.Range("N122").Value = -12151.5001

$ws.Range("H123").Value = 21529
$ws.Range("J123").Value = 21529
$ws.Range("L123").Value = 21529
$ws.Range("N123").Value = -26429

$ws.Range("H126").Value = 5109.1763
$ws.Range("I126").Value = 9122.846
$ws.Range("J126").Value = 2624.524
$ws.Range("K126").Value = 27368.538
$ws.Range("L126").Value = 7873.572
$ws.Range("M126").Value = -24898.538
$ws.Range("N126").Value = -12813.572

$ws.Range("H132").Value = 1800.3
$ws.Range("I132").Value = 1518.7354
$ws.Range("J132").Value = 2398.625
$ws.Range("K132").Value = 4556.206200000001
$ws.Range("L132").Value = 7195.875
$ws.Range("M132").Value = -2026.206200000001
$ws.Range("N132").Value = -12255.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 70010
$ws.Range("J5").Value = 70010
$ws.Range("L5").Value = 70010
$ws.Range("N5").Value = -70236

$ws.Range("H122").Value = 2471673.8
$ws.Range("I122").Value = 3110432
$ws.Range("J122").Value = 1002530
$ws.Range("K122").Value = 9331296
$ws.Range("L122").Value = 3007590
$ws.Range("M122").Value = -9328846
$ws.Range("N122").Value = -3012490

$ws.Range("H132").Value = 11372648
$ws.Range("I132").Value = 16196492
$ws.Range("J132").Value = 2157.5
$ws.Range("K132").Value = 48589476
$ws.Range("L132").Value = 6472.5
$ws.Range("M132").Value = -48586946
$ws.Range("N132").Value = -11532.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1154.1666
$ws.Range("I122").Value = 1150
$ws.Range("K122").Value = 3450
$ws.Range("M122").Value = -1000

$ws.Range("H132").Value = 1127.8536
$ws.Range("I132").Value = 863.8889
$ws.Range("J132").Value = 1636.9286
$ws.Range("K132").Value = 2591.6667
$ws.Range("L132").Value = 4910.7858
$ws.Range("M132").Value = -61.66670000000022
$ws.Range("N132").Value = -9970.7858
